$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44354
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 38000
$ws.Range("O2").Value = 38000
$ws.Range("P2").Value = 38000
$ws.Range("S2").Value = 2111

# Row 3
$ws.Range("D3").Value = 44435

# Row 4
$ws.Range("D4").Value = 44340
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = 37000
$ws.Range("O4").Value = 37000
$ws.Range("P4").Value = 37000
$ws.Range("S4").Value = 2056

# Row 5
$ws.Range("D5").Value = 44340
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("S5").Value = 1944

# Row 6
$ws.Range("D6").Value = 44333
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 38000
$ws.Range("O6").Value = 38000
$ws.Range("P6").Value = 38000
$ws.Range("S6").Value = 2111

# Row 7
$ws.Range("D7").Value = 44333
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 35000
$ws.Range("O7").Value = 35000
$ws.Range("P7").Value = 35000
$ws.Range("S7").Value = 1944

# Row 8
$ws.Range("D8").Value = 44445
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 32000
$ws.Range("O8").Value = 32000
$ws.Range("P8").Value = 32000
$ws.Range("S8").Value = 1778

# Row 9
$ws.Range("D9").Value = 44417
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 28000
$ws.Range("O9").Value = 28000
$ws.Range("P9").Value = 28000
$ws.Range("S9").Value = 1556

# Row 10
$ws.Range("D10").Value = 44389
$ws.Range("M10").Value = 35
$ws.Range("N10").Value = 29000
$ws.Range("O10").Value = 29000
$ws.Range("P10").Value = 29000
$ws.Range("S10").Value = 1611

# Row 11
$ws.Range("D11").Value = 44389
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 27000
$ws.Range("O11").Value = 27000
$ws.Range("P11").Value = 27000
$ws.Range("S11").Value = 1500

# Row 12
$ws.Range("D12").Value = 44410
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 15

# Row 13
$ws.Range("D13").Value = 44410
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 30000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 30000
$ws.Range("S13").Value = 1667

# Row 14
$ws.Range("D14").Value = 44410
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 28000
$ws.Range("O14").Value = 28000
$ws.Range("P14").Value = 28000
$ws.Range("S14").Value = 1556

# Row 17
$ws.Range("D17").Value = 44459
$ws.Range("M17").Value = 25

# Row 18
$ws.Range("D18").Value = 44382
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 20

# Row 19
$ws.Range("D19").Value = 44382
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 32000
$ws.Range("O19").Value = 32000
$ws.Range("P19").Value = 32000
$ws.Range("S19").Value = 1778

# Row 20
$ws.Range("D20").Value = 44382
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 30000
$ws.Range("O20").Value = 30000
$ws.Range("P20").Value = 30000
$ws.Range("S20").Value = 1667

# Row 21
$ws.Range("D21").Value = 44277
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 30000
$ws.Range("O21").Value = 30000
$ws.Range("P21").Value = 30000
$ws.Range("S21").Value = 1667

# Row 22
$ws.Range("D22").Value = 44277
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 28000
$ws.Range("O22").Value = 28000
$ws.Range("P22").Value = 28000
$ws.Range("S22").Value = 1556

# Row 23
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 16
$ws.Range("N23").Value = 35000
$ws.Range("O23").Value = 35000
$ws.Range("P23").Value = 35000
$ws.Range("S23").Value = 1944

# Row 24
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 25
$ws.Range("N24").Value = 30000
$ws.Range("O24").Value = 30000
$ws.Range("P24").Value = 30000
$ws.Range("S24").Value = 1667

# Row 25
$ws.Range("D25").Value = 44326
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 20
$ws.Range("N25").Value = 28000
$ws.Range("O25").Value = 28000
$ws.Range("P25").Value = 28000
$ws.Range("S25").Value = 1556

# Row 26
$ws.Range("D26").Value = 44319
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 140
$ws.Range("N26").Value = 27000
$ws.Range("O26").Value = 27000
$ws.Range("P26").Value = 27000
$ws.Range("S26").Value = 1500

# Row 27
$ws.Range("D27").Value = 44473
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 10
$ws.Range("N27").Value = 75000
$ws.Range("O27").Value = 75000
$ws.Range("P27").Value = 75000
$ws.Range("S27").Value = 4167

# Row 28
$ws.Range("D28").Value = 44473
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 25
$ws.Range("N28").Value = 72000
$ws.Range("O28").Value = 72000
$ws.Range("P28").Value = 72000
$ws.Range("S28").Value = 4000

# Row 29
$ws.Range("D29").Value = 44473
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 5
$ws.Range("N29").Value = 67000
$ws.Range("O29").Value = 67000
$ws.Range("P29").Value = 67000
$ws.Range("S29").Value = 3722

# Row 30
$ws.Range("D30").Value = 44452
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 36000
$ws.Range("O30").Value = 36000
$ws.Range("P30").Value = 36000
$ws.Range("S30").Value = 2000

# Row 31
$ws.Range("D31").Value = 44284
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 23000
$ws.Range("O31").Value = 23000
$ws.Range("P31").Value = 23000
$ws.Range("S31").Value = 1278

# Row 33
$ws.Range("D33").Value = 44312
$ws.Range("M33").Value = 160
$ws.Range("N33").Value = 26000
$ws.Range("O33").Value = 26000
$ws.Range("P33").Value = 26000
$ws.Range("S33").Value = 1444

# Row 34
$ws.Range("D34").Value = 44424
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 30
$ws.Range("N34").Value = 32000
$ws.Range("O34").Value = 32000
$ws.Range("P34").Value = 32000
$ws.Range("S34").Value = 1778

# Row 35
$ws.Range("D35").Value = 44396
$ws.Range("M35").Value = 35
$ws.Range("N35").Value = 37000
$ws.Range("O35").Value = 37000
$ws.Range("P35").Value = 37000
$ws.Range("S35").Value = 2056

# Row 36
$ws.Range("D36").Value = 44396
$ws.Range("N36").Value = 34000
$ws.Range("O36").Value = 34000
$ws.Range("P36").Value = 34000
$ws.Range("S36").Value = 1889

# Row 37
$ws.Range("D37").Value = 44431
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 32000
$ws.Range("O37").Value = 32000
$ws.Range("P37").Value = 32000
$ws.Range("S37").Value = 1778

# Row 38
$ws.Range("D38").Value = 44403
$ws.Range("L38").Value = "Especial"
$ws.Range("M38").Value = 25
$ws.Range("N38").Value = 33000
$ws.Range("O38").Value = 33000
$ws.Range("P38").Value = 33000
$ws.Range("S38").Value = 1833

# Row 39
$ws.Range("D39").Value = 44403
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 45

# Row 40
$ws.Range("D40").Value = 44403
$ws.Range("L40").Value = "Segunda"

# Row 44
$ws.Range("D44").Value = 44438
$ws.Range("M44").Value = 30
$ws.Range("N44").Value = 32000
$ws.Range("O44").Value = 32000
$ws.Range("P44").Value = 32000
$ws.Range("S44").Value = 1778

# Row 45
$ws.Range("D45").Value = 44529
$ws.Range("L45").Value = "Especial"
$ws.Range("M45").Value = 20
$ws.Range("N45").Value = 60000
$ws.Range("O45").Value = 60000
$ws.Range("P45").Value = 60000
$ws.Range("S45").Value = 3333

# Row 46
$ws.Range("D46").Value = 44529
$ws.Range("M46").Value = 50
$ws.Range("N46").Value = 58000
$ws.Range("O46").Value = 58000
$ws.Range("P46").Value = 58000
$ws.Range("S46").Value = 3222

# Row 47
$ws.Range("D47").Value = 44270
$ws.Range("L47").Value = "Especial"
$ws.Range("M47").Value = 70
$ws.Range("N47").Value = 38000
$ws.Range("O47").Value = 38000
$ws.Range("P47").Value = 38000
$ws.Range("S47").Value = 2111
